{"js": "// Replace the date line and every \"A\u00d7B=C\" answer cell in the practice\n// table with the new values from the target revision. Each old value is\n// unique in the document, so a simple exact-text search + replace per\n// pair is sufficient and avoids any row/column-index bookkeeping.\nconst replacements = [\n  [\"2025-01-17 Friday\", \"2025-01-18 Saturday\"],\n  [\"959\u00d75=4795\", \"322\u00d78=2576\"],\n  [\"129\u00d79=1161\", \"495\u00d74=1980\"],\n  [\"459\u00d78=3672\", \"358\u00d76=2148\"],\n  [\"944\u00d79=8496\", \"381\u00d78=3048\"],\n  [\"565\u00d74=2260\", \"789\u00d73=2367\"],\n  [\"718\u00d76=4308\", \"738\u00d74=2952\"],\n  [\"437\u00d73=1311\", \"156\u00d72=312\"],\n  [\"526\u00d78=4208\", \"426\u00d75=2130\"],\n  [\"327\u00d75=1635\", \"587\u00d72=1174\"],\n  [\"240\u00d75=1200\", \"140\u00d72=280\"],\n  [\"231\u00d75=1155\", \"516\u00d72=1032\"],\n  [\"337\u00d72=674\", \"485\u00d73=1455\"],\n  [\"707\u00d79=6363\", \"848\u00d74=3392\"],\n  [\"651\u00d74=2604\", \"259\u00d79=2331\"],\n  [\"389\u00d76=2334\", \"630\u00d79=5670\"],\n  [\"435\u00d74=1740\", \"529\u00d72=1058\"],\n  [\"696\u00d79=6264\", \"451\u00d76=2706\"],\n  [\"337\u00d74=1348\", \"540\u00d74=2160\"],\n  [\"883\u00d75=4415\", \"681\u00d73=2043\"],\n  [\"238\u00d73=714\", \"774\u00d78=6192\"],\n  [\"855\u00d74=3420\", \"487\u00d77=3409\"],\n  [\"799\u00d75=3995\", \"390\u00d78=3120\"],\n  [\"692\u00d77=4844\", \"972\u00d74=3888\"],\n  [\"224\u00d79=2016\", \"949\u00d75=4745\"],\n  [\"754\u00d74=3016\", \"496\u00d77=3472\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A\u00d7B=C\" answer cell in the practice\n# table with the new values from the target revision. Each old value is\n# unique in the document, so a simple Find/Replace (wdReplaceAll) per\n# pair is sufficient and avoids any row/column-index bookkeeping.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-01-17 Friday\", \"2025-01-18 Saturday\"),\n    @(\"959\u00d75=4795\", \"322\u00d78=2576\"),\n    @(\"129\u00d79=1161\", \"495\u00d74=1980\"),\n    @(\"459\u00d78=3672\", \"358\u00d76=2148\"),\n    @(\"944\u00d79=8496\", \"381\u00d78=3048\"),\n    @(\"565\u00d74=2260\", \"789\u00d73=2367\"),\n    @(\"718\u00d76=4308\", \"738\u00d74=2952\"),\n    @(\"437\u00d73=1311\", \"156\u00d72=312\"),\n    @(\"526\u00d78=4208\", \"426\u00d75=2130\"),\n    @(\"327\u00d75=1635\", \"587\u00d72=1174\"),\n    @(\"240\u00d75=1200\", \"140\u00d72=280\"),\n    @(\"231\u00d75=1155\", \"516\u00d72=1032\"),\n    @(\"337\u00d72=674\", \"485\u00d73=1455\"),\n    @(\"707\u00d79=6363\", \"848\u00d74=3392\"),\n    @(\"651\u00d74=2604\", \"259\u00d79=2331\"),\n    @(\"389\u00d76=2334\", \"630\u00d79=5670\"),\n    @(\"435\u00d74=1740\", \"529\u00d72=1058\"),\n    @(\"696\u00d79=6264\", \"451\u00d76=2706\"),\n    @(\"337\u00d74=1348\", \"540\u00d74=2160\"),\n    @(\"883\u00d75=4415\", \"681\u00d73=2043\"),\n    @(\"238\u00d73=714\", \"774\u00d78=6192\"),\n    @(\"855\u00d74=3420\", \"487\u00d77=3409\"),\n    @(\"799\u00d75=3995\", \"390\u00d78=3120\"),\n    @(\"692\u00d77=4844\", \"972\u00d74=3888\"),\n    @(\"224\u00d79=2016\", \"949\u00d75=4745\"),\n    @(\"754\u00d74=3016\", \"496\u00d77=3472\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
